$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Samiksha Pansare) shares every column from E through T with the
# new intern row being added, so use it as a starting template: copy it
# down to row 4, then overwrite the columns that actually differ (A-D).
$ws.Range("A3:T3").Copy($ws.Range("A4:T4"))

$ws.Range("A4").Value = "Deep"
$ws.Range("B4").Value = "19102B0052"

# The mobile number has a leading zero, which must be kept literal (as
# text) rather than collapsed by numeric auto-conversion. Build the text
# value with TEXT() in a scratch cell, then paste only the resulting
# value into C4 so no stray cell/number formatting is left behind.
$ws.Range("ZZ1").Formula = '=TEXT(9323165165,"00000000000")'
$ws.Range("ZZ1").Copy()
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()

$ws.Range("D4").Value = "esotericdeep@gmail.com"
